# Update "want to go" counts (column F) across sheets to reflect the
# regenerated data output (commit: "Update gh-pages to output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 2615
$ws1.Range("F6").Value  = 191
$ws1.Range("F13").Value = 5649
$ws1.Range("F15").Value = 1742
$ws1.Range("F16").Value = 4092
$ws1.Range("F20").Value = 4738
$ws1.Range("F21").Value = 6152
$ws1.Range("F23").Value = 1049
$ws1.Range("F26").Value = 492
$ws1.Range("F36").Value = 1695
$ws1.Range("F42").Value = 624
$ws1.Range("F44").Value = 3362
$ws1.Range("F46").Value = 278
$ws1.Range("F48").Value = 7
$ws1.Range("F49").Value = 3876

# --- Sheet "本地生活" (Local life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 3805

# --- Sheet "全部类型" (All types, merged listing) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 3805
$ws4.Range("F4").Value  = 2615
$ws4.Range("F10").Value = 191
$ws4.Range("F18").Value = 1742
$ws4.Range("F19").Value = 4738
$ws4.Range("F21").Value = 1049
$ws4.Range("F24").Value = 492
$ws4.Range("F35").Value = 1695
$ws4.Range("F39").Value = 624
$ws4.Range("F43").Value = 3362
$ws4.Range("F46").Value = 278
$ws4.Range("F49").Value = 3876

$wb.Save()
